# Implement AEBC "Alternate Economic Base Case" control lever (#51)
#
# The "Key to Variables" sheet documents every acronym/variable used by the
# model. A new plcy-ctrl-ctr (control lever) row needs to be added for the
# new "AEBC" / "Alternate Economic Base Case" lever, directly above the
# existing BAEPAbCiGC row (which was previously the first plcy-ctrl-ctr
# entry, at row 142).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Key to Variables")

# Insert a new blank row at 142; this shifts row 142 (BAEPAbCiGC) and every
# row below it down by one (142->143, ..., 194->195).
$ws.Rows.Item(142).Insert()

# Column layout (row 1 headers): A=Top Level Folder, B=Acronym,
# C=Variable/Lever name, D/E=extra notes, F=Level, G=Description.
$ws.Cells.Item(142, 1).Value = "plcy-ctrl-ctr"
$ws.Cells.Item(142, 2).Value = "AEBC"
$ws.Cells.Item(142, 3).Value = "Alternate Economic Base Case"
$ws.Cells.Item(142, 6).Value = "very high"

# The "Level" column is color coded; "very high" entries use a red fill.
$ws.Cells.Item(142, 6).Interior.Color = 255

# Leave the "About" sheet as the selected tab.
$about = $wb.Worksheets.Item("About")
$about.Activate()
